$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("85").Insert()

$ws.Range("A85").Value = 8
$ws.Range("B85").Value = "Terminal La Palmera de La Serena"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value = 44644
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100109
$ws.Range("H85").Value = "Uva"
$ws.Range("I85").Value = 100109001
$ws.Range("J85").Value = "Uva"
$ws.Range("K85").Value = "Red Globe"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 360
$ws.Range("N85").Value = 9500
$ws.Range("O85").Value = 10000
$ws.Range("P85").Value = 9750
$ws.Range("Q85").Value = "`$/bandeja 18 kilos"
$ws.Range("R85").Value = "Provincia del Elquí"
$ws.Range("S85").Value = 542
$ws.Range("T85").Value = 18
